$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions data pull).
# Price (D) and Volume(1h) (E) are plain-text columns in this sheet, so a
# leading apostrophe is used on every assignment (just like typing '1.041
# into a cell) to stop Excel from auto-coercing numeric-looking strings
# (e.g. "1.041", "0.000009150", "40.00") into Number/Date values, which
# would silently drop the exact text formatting these cells rely on.

$ws.Range("D2").Value = "'27.794.36"
$ws.Range("E2").Value = "'  +3.27%  "

$ws.Range("D3").Value = "'1.868.12"
$ws.Range("E3").Value = "'  +3.02%  "

$ws.Range("D4").Value = "'1.041"
$ws.Range("E4").Value = "'  +3.06%  "

$ws.Range("D5").Value = "'325.49"
$ws.Range("E5").Value = "'  +4.39%  "

$ws.Range("D6").Value = "'1.037"
$ws.Range("E6").Value = "'  +3.11%  "

$ws.Range("D7").Value = "'0.4424"
$ws.Range("E7").Value = "'  +3.00%  "

$ws.Range("D8").Value = "'0.3803"
$ws.Range("E8").Value = "'  +2.89%  "

$ws.Range("D9").Value = "'0.07470"
$ws.Range("E9").Value = "'  +2.96%  "

$ws.Range("D10").Value = "'0.8865"
$ws.Range("E10").Value = "'  +2.28%  "

$ws.Range("D11").Value = "'21.82"
$ws.Range("E11").Value = "'  +2.44%  "

$ws.Range("D12").Value = "'1.870.68"
$ws.Range("E12").Value = "'  -12.71%  "

$ws.Range("D13").Value = "'5.568"
$ws.Range("E13").Value = "'  +2.98%  "

$ws.Range("D14").Value = "'6.761"
$ws.Range("E14").Value = "'  +1.77%  "

$ws.Range("D15").Value = "'0.07244"
$ws.Range("E15").Value = "'  +4.02%  "

$ws.Range("D16").Value = "'83.86"
$ws.Range("E16").Value = "'  +3.69%  "

$ws.Range("D17").Value = "'1.041"
$ws.Range("E17").Value = "'  +3.42%  "

$ws.Range("D18").Value = "'0.000009150"
$ws.Range("E18").Value = "'  +2.85%  "

$ws.Range("D19").Value = "'1.037"
$ws.Range("E19").Value = "'  +3.16%  "

$ws.Range("D20").Value = "'15.57"
$ws.Range("E20").Value = "'  +1.91%  "

$ws.Range("D21").Value = "'27.817.48"
$ws.Range("E21").Value = "'  +3.19%  "

$ws.Range("D22").Value = "'5.325"
$ws.Range("E22").Value = "'  +2.56%  "

$ws.Range("D23").Value = "'11.44"
$ws.Range("E23").Value = "'  +4.13%  "

$ws.Range("D24").Value = "'1.966"
$ws.Range("E24").Value = "'  +4.15%  "

$ws.Range("D25").Value = "'158.70"
$ws.Range("E25").Value = "'  +2.82%  "

$ws.Range("D26").Value = "'18.90"
$ws.Range("E26").Value = "'  +2.83%  "

$ws.Range("D27").Value = "'1.994"
$ws.Range("E27").Value = "'  +2.59%  "

$ws.Range("D28").Value = "'5.332"
$ws.Range("E28").Value = "'  +1.92%  "

$ws.Range("D29").Value = "'117.73"
$ws.Range("E29").Value = "'  +2.54%  "

$ws.Range("D30").Value = "'0.09122"
$ws.Range("E30").Value = "'  +1.82%  "

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'0.7768"
$ws.Range("E31").Value = "'  +4.50%  "

$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = "'1.217"
$ws.Range("E32").Value = "'  +4.37%  "

$ws.Range("D33").Value = "'3.038"
$ws.Range("E33").Value = "'  +8.40%  "

$ws.Range("D34").Value = "'4.598"
$ws.Range("E34").Value = "'  +3.70%  "

$ws.Range("D35").Value = "'1.038"
$ws.Range("E35").Value = "'  +3.28%  "

$ws.Range("D36").Value = "'1.170"
$ws.Range("E36").Value = "'  +4.41%  "

$ws.Range("D37").Value = "'0.01996"
$ws.Range("E37").Value = "'  +3.62%  "

$ws.Range("D38").Value = "'0.05367"
$ws.Range("E38").Value = "'  +2.52%  "

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = "'0.5215"
$ws.Range("E39").Value = "'  +2.26%  "

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = "'2.845"
$ws.Range("E40").Value = "'  +2.82%  "

$ws.Range("D41").Value = "'0.1696"
$ws.Range("E41").Value = "'  +2.67%  "

$ws.Range("D42").Value = "'6.894"
$ws.Range("E42").Value = "'  +6.50%  "

$ws.Range("D43").Value = "'8.734"
$ws.Range("E43").Value = "'  +4.97%  "

$ws.Range("D44").Value = "'109.89"
$ws.Range("E44").Value = "'  +2.46%  "

$ws.Range("D45").Value = "'10.63"
$ws.Range("E45").Value = "'  +1.35%  "

$ws.Range("D46").Value = "'1.727"
$ws.Range("E46").Value = "'  +4.76%  "

$ws.Range("D47").Value = "'0.4716"
$ws.Range("E47").Value = "'  +2.87%  "

$ws.Range("D48").Value = "'0.06447"
$ws.Range("E48").Value = "'  +2.43%  "

$ws.Range("D49").Value = "'1.889"
$ws.Range("E49").Value = "'  +4.77%  "

$ws.Range("D50").Value = "'40.00"
$ws.Range("E50").Value = "'  +5.21%  "

$ws.Range("D51").Value = "'64.67"
$ws.Range("E51").Value = "'  +1.56%  "
